$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H100").Value = 7335
$ws.Range("I100").Value = 6002.5
$ws.Range("J100").Value = 10000
$ws.Range("K100").Value = 6002.5
$ws.Range("L100").Value = 10000
$ws.Range("M100").Value = -5461.5
$ws.Range("N100").Value = -11082

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H129").Value = 743.2
$ws.Range("J129").Value = 987
$ws.Range("L129").Value = 2961
$ws.Range("N129").Value = -12961

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 14044.362
$ws.Range("I132").Value = 15181.137
$ws.Range("J132").Value = 2189.4285
$ws.Range("K132").Value = 45543.411
$ws.Range("L132").Value = 6568.2855
$ws.Range("M132").Value = -43013.411
$ws.Range("N132").Value = -11628.2855

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 5893.6943
$ws.Range("I138").Value = 4339.8184
$ws.Range("J138").Value = 6173.902
$ws.Range("K138").Value = 13019.4552
$ws.Range("L138").Value = 18521.706
$ws.Range("M138").Value = -7879.4552
$ws.Range("N138").Value = -28801.706

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H140").Value = 64140
$ws.Range("J140").Value = 64140
$ws.Range("L140").Value = 64140
$ws.Range("N140").Value = -74500

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1265.8889
$ws.Range("I2").Value = 1434.0714
$ws.Range("J2").Value = 677.25
$ws.Range("K2").Value = 1434.0714
$ws.Range("L2").Value = 677.25
$ws.Range("M2").Value = -1321.0714
$ws.Range("N2").Value = -903.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 12642.053
$ws.Range("I32").Value = 9082.655000000001
$ws.Range("J32").Value = 24111.223
$ws.Range("K32").Value = 9082.655000000001
$ws.Range("L32").Value = 24111.223
$ws.Range("M32").Value = -8795.655000000001
$ws.Range("N32").Value = -24685.223

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H52").Value = 41796.668
$ws.Range("J52").Value = 41796.668
$ws.Range("L52").Value = 41796.668
$ws.Range("N52").Value = -42432.668

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 1462.174
$ws.Range("I74").Value = 949.5714
$ws.Range("J74").Value = 3670.3076
$ws.Range("K74").Value = 949.5714
$ws.Range("L74").Value = 3670.3076
$ws.Range("M74").Value = -75.57140000000004
$ws.Range("N74").Value = -5418.3076

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 1462.174
$ws.Range("I77").Value = 949.5714
$ws.Range("J77").Value = 3670.3076
$ws.Range("K77").Value = 4747.857
$ws.Range("L77").Value = 18351.538
$ws.Range("M77").Value = -379.857
$ws.Range("N77").Value = -27087.538

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H109").Value = 15377
$ws.Range("J109").Value = 15377
$ws.Range("L109").Value = 15377
$ws.Range("N109").Value = -18151

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H116").Value = 1265.8889
$ws.Range("I116").Value = 1434.0714
$ws.Range("J116").Value = 677.25
$ws.Range("K116").Value = 1434.0714
$ws.Range("L116").Value = 677.25
$ws.Range("M116").Value = 859.9286
$ws.Range("N116").Value = -5265.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 2731.8462
$ws.Range("I122").Value = 1944.4445
$ws.Range("K122").Value = 5833.333500000001
$ws.Range("M122").Value = -3383.333500000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 2470.1843
$ws.Range("I132").Value = 1969.2084
$ws.Range("J132").Value = 3329
$ws.Range("K132").Value = 5907.6252
$ws.Range("L132").Value = 9987
$ws.Range("M132").Value = -3377.6252
$ws.Range("N132").Value = -15047

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1265.8889
$ws.Range("I3").Value = 1434.0714
$ws.Range("J3").Value = 677.25
$ws.Range("K3").Value = 1434.0714
$ws.Range("L3").Value = 677.25
$ws.Range("M3").Value = -1320.0714
$ws.Range("N3").Value = -905.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 2164.7334
$ws.Range("I99").Value = 1312.5
$ws.Range("J99").Value = 2474.6365
$ws.Range("K99").Value = 1312.5
$ws.Range("L99").Value = 2474.6365
$ws.Range("M99").Value = 185.5
$ws.Range("N99").Value = -5470.636500000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 700
$ws.Range("I107").Value = 700
$ws.Range("K107").Value = 700
$ws.Range("M107").Value = 1220

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2372.2354
$ws.Range("I134").Value = 1836.4445
$ws.Range("J134").Value = 2975
$ws.Range("K134").Value = 5509.333500000001
$ws.Range("L134").Value = 8925
$ws.Range("M134").Value = -2974.333500000001
$ws.Range("N134").Value = -13995

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H140").Value = 59765
$ws.Range("J140").Value = 59765
$ws.Range("L140").Value = 59765
$ws.Range("N140").Value = -70125

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H141").Value = 42736.25
$ws.Range("J141").Value = 42736.25
$ws.Range("L141").Value = 42736.25
$ws.Range("N141").Value = -53096.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5422.6304
$ws.Range("I31").Value = 2200.5186
$ws.Range("J31").Value = 10001.421
$ws.Range("K31").Value = 2200.5186
$ws.Range("L31").Value = 10001.421
$ws.Range("M31").Value = -1905.5186
$ws.Range("N31").Value = -10591.421

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 5422.6304
$ws.Range("I34").Value = 2200.5186
$ws.Range("J34").Value = 10001.421
$ws.Range("K34").Value = 2200.5186
$ws.Range("L34").Value = 10001.421
$ws.Range("M34").Value = -1998.5186
$ws.Range("N34").Value = -10405.421

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 9251.666999999999
$ws.Range("I99").Value = 2272.9167
$ws.Range("J99").Value = 37166.668
$ws.Range("K99").Value = 2272.9167
$ws.Range("L99").Value = 37166.668
$ws.Range("M99").Value = -774.9167000000002
$ws.Range("N99").Value = -40162.668

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 635.3333
$ws.Range("I107").Value = 360.83334
$ws.Range("J107").Value = 1733.3334
$ws.Range("K107").Value = 360.83334
$ws.Range("L107").Value = 1733.3334
$ws.Range("M107").Value = 1559.16666
$ws.Range("N107").Value = -5573.3334

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value = 9251.666999999999
$ws.Range("I126").Value = 2272.9167
$ws.Range("J126").Value = 37166.668
$ws.Range("K126").Value = 6818.750100000001
$ws.Range("L126").Value = 111500.004
$ws.Range("M126").Value = -4348.750100000001
$ws.Range("N126").Value = -116440.004

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H133").Value = 25299.666
$ws.Range("J133").Value = 25299.666
$ws.Range("L133").Value = 25299.666
$ws.Range("N133").Value = -30359.666

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 2740.8518
$ws.Range("I134").Value = 2791.7917
$ws.Range("J134").Value = 2333.3333
$ws.Range("K134").Value = 8375.375100000001
$ws.Range("L134").Value = 6999.999899999999
$ws.Range("M134").Value = -5840.375100000001
$ws.Range("N134").Value = -12069.9999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 91.86207
$ws.Range("I2").Value = 269.1111
$ws.Range("J2").Value = 12.1
$ws.Range("K2").Value = 1614.6666
$ws.Range("L2").Value = 72.59999999999999
$ws.Range("M2").Value = -1501.6666
$ws.Range("N2").Value = -298.6

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H86").Value = 942.8095
$ws.Range("I86").Value = 716.6667
$ws.Range("J86").Value = 1033.2667
$ws.Range("K86").Value = 2150.0001
$ws.Range("L86").Value = 3099.800099999999
$ws.Range("M86").Value = -964.0001000000002
$ws.Range("N86").Value = -5471.800099999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H89").Value = 942.8095
$ws.Range("I89").Value = 716.6667
$ws.Range("J89").Value = 1033.2667
$ws.Range("K89").Value = 6450.0003
$ws.Range("L89").Value = 9299.400299999999
$ws.Range("M89").Value = -522.0002999999997
$ws.Range("N89").Value = -21155.4003

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H97").Value = 2085.25
$ws.Range("I97").Value = 1590
$ws.Range("J97").Value = 2184.3
$ws.Range("K97").Value = 4770
$ws.Range("L97").Value = 6552.900000000001
$ws.Range("M97").Value = -4274
$ws.Range("N97").Value = -7544.900000000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H137").Value = 6575.905
$ws.Range("I137").Value = 2127.0908
$ws.Range("J137").Value = 8154.516
$ws.Range("K137").Value = 6381.2724
$ws.Range("L137").Value = 24463.548
$ws.Range("M137").Value = -1281.2724
$ws.Range("N137").Value = -34663.548

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H14").Value = 1467185
$ws.Range("I14").Value = 3333431.8
$ws.Range("J14").Value = 67500
$ws.Range("K14").Value = 3333431.8
$ws.Range("L14").Value = 67500
$ws.Range("M14").Value = -3333263.8
$ws.Range("N14").Value = -67836

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2066.6667
$ws.Range("I80").Value = 2000
$ws.Range("J80").Value = 2080
$ws.Range("K80").Value = 2000
$ws.Range("L80").Value = 2080
$ws.Range("M80").Value = -1002
$ws.Range("N80").Value = -4076

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value = 2066.6667
$ws.Range("I83").Value = 2000
$ws.Range("J83").Value = 2080
$ws.Range("K83").Value = 10000
$ws.Range("L83").Value = 10400
$ws.Range("M83").Value = -5008
$ws.Range("N83").Value = -20384

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 2365.0571
$ws.Range("I122").Value = 1704.8182
$ws.Range("K122").Value = 5114.4546
$ws.Range("M122").Value = -2664.4546

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H124").Value = 38980
$ws.Range("J124").Value = 38980
$ws.Range("L124").Value = 38980
$ws.Range("N124").Value = -48800

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2554.6738
$ws.Range("I132").Value = 2259.7666
$ws.Range("J132").Value = 3107.625
$ws.Range("K132").Value = 6779.2998
$ws.Range("L132").Value = 9322.875
$ws.Range("M132").Value = -4249.2998
$ws.Range("N132").Value = -14382.875

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H127").Value = 39200
$ws.Range("J127").Value = 39200
$ws.Range("L127").Value = 39200
$ws.Range("N127").Value = -49120

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1746.6061
$ws.Range("I132").Value = 1225.16
$ws.Range("J132").Value = 3376.125
$ws.Range("K132").Value = 3675.48
$ws.Range("L132").Value = 10128.375
$ws.Range("M132").Value = -1145.48
$ws.Range("N132").Value = -15188.375
